# Add a "Full size /w index" column (new column C) with timing data,
# shifting the existing "Half size" / "Quarter size " columns one to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C (Half size). This shifts
# the old C -> D and old D -> E, carries the merged header cell along, and
# updates the sheet dimension automatically.
$ws.Columns("C").Insert()

# Match the target column width for the newly inserted column (xml width 15.5).
$ws.Columns("C").ColumnWidth = 14.666666666666666

# Header cell for the new column (row 2).
$ws.Range("C2").Value = "Full size /w index"

# New column C data values (rows 3-16).
$ws.Range("C3").Value = 3125.2950000000001
$ws.Range("C4").Value = 230.09100000000001
$ws.Range("C5").Value = 94.924000000000007
$ws.Range("C6").Value = 5.1360000000000001
$ws.Range("C7").Value = 644.971
$ws.Range("C8").Value = 3008.6819999999998
$ws.Range("C9").Value = 11585.268
$ws.Range("C10").Value = 3373.23
$ws.Range("C11").Value = 26880.744999999999
$ws.Range("C12").Value = 1646.471
$ws.Range("C13").Value = 435.64699999999999
$ws.Range("C14").Value = 16830.469000000001
$ws.Range("C15").Value = 1539.2619999999999
$ws.Range("C16").Value = 518.20299999999997

# Move the selection to E16, matching the committed selection state.
[void]$ws.Range("E16").Select()
